$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "eyes_closed_rest"
$ws.Range("B14").Value = "Eyes_Closed_Rest"
$ws.Range("C14").Value = "None"

$ws.Range("C14").Select()
